# Actualización automática 2025-10-29 09:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("L4").Value = 1187.48
$wsVentasGrupo.Range("I31").Value = 218.62
$wsVentasGrupo.Range("L31").Value = 447.78
$wsVentasGrupo.Range("I60").Value = "3 de 58"
$wsVentasGrupo.Range("L60").Value = "9 de 58"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F4").Value = 1270.98
$wsVentaMensual.Range("F31").Value = 8844.76
$wsVentaMensual.Range("F60").Value = 55720.09

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D7").Value = 1022.29
$wsCumplimiento.Range("E7").Value = -135.578983712426
$wsCumplimiento.Range("F7").Value = 1.152900980389371

$wsCumplimiento.Range("D11").Value = 13461.01
$wsCumplimiento.Range("E11").Value = 6112.050249249698
$wsCumplimiento.Range("F11").Value = 0.6877314956671635

$wsCumplimiento.Range("D12").Value = 26412.95
$wsCumplimiento.Range("E12").Value = 22211.11
$wsCumplimiento.Range("F12").Value = 0.5432074162461958

$wsCumplimiento.Range("D14").Value = 61324.25000000001
$wsCumplimiento.Range("E14").Value = 38573.74284188786
$wsCumplimiento.Range("F14").Value = 0.6138686900052147
